$d = $word.ActiveDocument
$s = $d.Shapes.Item(33)
$tr = $s.TextFrame.TextRange
Write-Output "Text=[$($tr.Text)]"
Write-Output "Start=$($tr.Start) End=$($tr.End)"
$tr2 = $s.TextFrame.TextRange
Write-Output "Text2=[$($tr2.Text)]"
